$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ackley")
$ws.Range("G2").Value = 0.07596752595681984
$ws.Range("H2").Value = 21.7929
$ws.Range("I2").Value = 21.7227
$ws.Range("G3").Value = 0.09212023664754337
$ws.Range("H3").Value = 21.71274
$ws.Range("I3").Value = 21.5946
$ws.Range("G4").Value = 0.1586123355858529
$ws.Range("H4").Value = 21.72526
$ws.Range("I4").Value = 21.5619
$ws.Range("G5").Value = 0.1533800345546995
$ws.Range("H5").Value = 21.7266
$ws.Range("I5").Value = 21.4952
$ws.Range("G6").Value = 0.08188291641119327
$ws.Range("H6").Value = 21.80532
$ws.Range("I6").Value = 21.7122
$ws.Range("G7").Value = 0.1116575702762732
$ws.Range("H7").Value = 21.71374
$ws.Range("I7").Value = 21.5609

$ws = $wb.Worksheets.Item("damavandi")
$ws.Range("G2").Value = 17.93712769196897
$ws.Range("H2").Value = 39.76456
$ws.Range("I2").Value = 16.018
$ws.Range("G3").Value = 20.85593328539867
$ws.Range("H3").Value = 53.2081
$ws.Range("I3").Value = 19.8783
$ws.Range("G4").Value = 28.90766451052385
$ws.Range("H4").Value = 31.61004
$ws.Range("I4").Value = 7.9318
$ws.Range("G5").Value = 37.82769697582712
$ws.Range("H5").Value = 55.25259999999999
$ws.Range("I5").Value = 4.4936
$ws.Range("G6").Value = 30.24290975291564
$ws.Range("H6").Value = 48.01456
$ws.Range("I6").Value = 18.1499
$ws.Range("G7").Value = 21.1288959196878
$ws.Range("H7").Value = 33.0321
$ws.Range("I7").Value = 3.5803

$ws = $wb.Worksheets.Item("griewank")
$ws.Range("G2").Value = 87.92449984035734
$ws.Range("H2").Value = 1022.16538
$ws.Range("I2").Value = 883.2292
$ws.Range("G3").Value = 140.9240395389055
$ws.Range("H3").Value = 917.2071600000002
$ws.Range("I3").Value = 727.1133
$ws.Range("G4").Value = 104.6777782946313
$ws.Range("H4").Value = 901.1131000000001
$ws.Range("I4").Value = 785.1153
$ws.Range("G5").Value = 102.8395656484654
$ws.Range("H5").Value = 961.5845999999999
$ws.Range("I5").Value = 858.9081
$ws.Range("G6").Value = 142.2109641956029
$ws.Range("H6").Value = 990.3490400000001
$ws.Range("I6").Value = 871.9924
$ws.Range("G7").Value = 284.0696350939763
$ws.Range("H7").Value = 781.6322
$ws.Range("I7").Value = 433.283

$ws = $wb.Worksheets.Item("schwefel")
$ws.Range("G2").Value = 877.1541750300074
$ws.Range("H2").Value = 13128.63494
$ws.Range("I2").Value = 12256.6829
$ws.Range("G3").Value = 1320.17973868919
$ws.Range("H3").Value = 12040.1156
$ws.Range("I3").Value = 10945.053
$ws.Range("G4").Value = 481.083244658903
$ws.Range("H4").Value = 12678.81668
$ws.Range("I4").Value = 11919.8706
$ws.Range("G5").Value = 1106.947020254855
$ws.Range("H5").Value = 12628.67994
$ws.Range("I5").Value = 11573.1087
$ws.Range("G6").Value = 1637.142369517997
$ws.Range("H6").Value = 12216.30614
$ws.Range("I6").Value = 10390.432
$ws.Range("G7").Value = 867.3822948542314
$ws.Range("H7").Value = 12641.17132
$ws.Range("I7").Value = 11655.1163

$ws = $wb.Worksheets.Item("rastrigin")
$ws.Range("G2").Value = 40.87643723742322
$ws.Range("H2").Value = 507.6675999999999
$ws.Range("I2").Value = 472.4385
$ws.Range("G3").Value = 65.21131537385213
$ws.Range("H3").Value = 468.2537600000001
$ws.Range("I3").Value = 386.896
$ws.Range("G4").Value = 41.33860110935058
$ws.Range("H4").Value = 552.9428399999999
$ws.Range("I4").Value = 511.4706
$ws.Range("G5").Value = 46.52777561115509
$ws.Range("H5").Value = 558.09998
$ws.Range("I5").Value = 510.847
$ws.Range("G6").Value = 28.21055818453442
$ws.Range("H6").Value = 522.77266
$ws.Range("I6").Value = 490.3805
$ws.Range("G7").Value = 61.38756317278772
$ws.Range("H7").Value = 542.97506
$ws.Range("I7").Value = 484.4864

$ws = $wb.Worksheets.Item("sphere")
$ws.Range("G2").Value = 25.89610202766045
$ws.Range("H2").Value = 233.36342
$ws.Range("I2").Value = 198.8525
$ws.Range("G3").Value = 63.49398550364436
$ws.Range("H3").Value = 263.88612
$ws.Range("I3").Value = 207.3536
$ws.Range("G4").Value = 22.85236975155978
$ws.Range("H4").Value = 226.64798
$ws.Range("I4").Value = 209.1884
$ws.Range("G5").Value = 38.1595477544874
$ws.Range("H5").Value = 191.83258
$ws.Range("I5").Value = 132.337
$ws.Range("G6").Value = 26.57540943330888
$ws.Range("H6").Value = 234.69674
$ws.Range("I6").Value = 205.4567
$ws.Range("G7").Value = 16.19556663185955
$ws.Range("H7").Value = 235.92808
$ws.Range("I7").Value = 209.0612

$ws = $wb.Worksheets.Item("rotatedhyperellipsoid")
$ws.Range("G2").Value = 99425.41543711444
$ws.Range("H2").Value = 633729.06134
$ws.Range("I2").Value = 480832.4097
$ws.Range("G3").Value = 52818.14373243946
$ws.Range("H3").Value = 743464.65056
$ws.Range("I3").Value = 675409.7943
$ws.Range("G4").Value = 129150.6663096203
$ws.Range("H4").Value = 617044.5622400001
$ws.Range("I4").Value = 420118.3565
$ws.Range("G5").Value = 221867.6249859901
$ws.Range("H5").Value = 668533.86226
$ws.Range("I5").Value = 490862.3612
$ws.Range("G6").Value = 86551.05731756629
$ws.Range("H6").Value = 676537.74808
$ws.Range("I6").Value = 532093.8329
$ws.Range("G7").Value = 84024.32754457422
$ws.Range("H7").Value = 662551.3906400001
$ws.Range("I7").Value = 538511.7928000001

$ws = $wb.Worksheets.Item("perm")
$ws.Range("G2").Value = 0.5971260687660521
$ws.Range("H2").Value = 0.69272
$ws.Range("I2").Value = 0.1839
$ws.Range("G3").Value = 0.8458535168692036
$ws.Range("H3").Value = 0.68302
$ws.Range("I3").Value = 0.0046
$ws.Range("G4").Value = 0.2342844147612047
$ws.Range("H4").Value = 0.2928799999999999
$ws.Range("I4").Value = 0.0486
$ws.Range("G5").Value = 0.7047567573283707
$ws.Range("H5").Value = 0.68652
$ws.Range("I5").Value = 0.109
$ws.Range("G6").Value = 0.3872438095050715
$ws.Range("H6").Value = 0.43524
$ws.Range("I6").Value = 0.0071
$ws.Range("G7").Value = 0.4489049710127969
$ws.Range("H7").Value = 0.38404
$ws.Range("I7").Value = 0.0081

$ws = $wb.Worksheets.Item("zakharov")
$ws.Range("G2").Value = 99259477924.72134
$ws.Range("H2").Value = 117523181459.0076
$ws.Range("I2").Value = 3258581024.9093
$ws.Range("G3").Value = 130257738920.7657
$ws.Range("H3").Value = 212833133361.048
$ws.Range("I3").Value = 22301875699.1304
$ws.Range("G4").Value = 49780609578.71429
$ws.Range("H4").Value = 122589279656.7849
$ws.Range("I4").Value = 48541632039.7682
$ws.Range("G5").Value = 129150725140.3167
$ws.Range("H5").Value = 107327363756.9536
$ws.Range("I5").Value = 12077815403.9756
$ws.Range("G6").Value = 102958172982.9337
$ws.Range("H6").Value = 116580862034.8071
$ws.Range("I6").Value = 6524159381.6133
$ws.Range("G7").Value = 78115976885.33308
$ws.Range("H7").Value = 68710646755.63744
$ws.Range("I7").Value = 2726514617.1257

$ws = $wb.Worksheets.Item("rosenbrock")
$ws.Range("G2").Value = 955023623225895.8
$ws.Range("H2").Value = 9861233160861486
$ws.Range("I2").Value = 8729391134443771
$ws.Range("G3").Value = 3012134138115961
$ws.Range("H3").Value = 9582163885512764
$ws.Range("I3").Value = 6820086315159471
$ws.Range("G4").Value = 1504211702602553
$ws.Range("H4").Value = 8514079842554832
$ws.Range("I4").Value = 6086606790638312
$ws.Range("G5").Value = 3164412425455727
$ws.Range("H5").Value = 11243938308501220
$ws.Range("I5").Value = 7580086636739949
$ws.Range("G6").Value = 1477010380859709
$ws.Range("H6").Value = 8575000436035427
$ws.Range("I6").Value = 7009619868328849
$ws.Range("G7").Value = 2592002394454354
$ws.Range("H7").Value = 10600695485231290
$ws.Range("I7").Value = 7338958942283823
